$d = $word.ActiveDocument
$c = $d.Content

# ---------------------------------------------------------------------------
# Cosmetic clean-up: a handful of sentences elsewhere in the document had
# their text split across extra runs wrapped in <w:proofErr/> spell/grammar
# -check bookmarks (e.g. around "Moodle", "runServer", "runManagerSwing",
# "ResultSet"/"ResultSetMetaData", and the "entire table contents is"
# grammar flag). None of these change the visible wording, so re-asserting
# the same text via Find/Replace simply collapses the runs back together
# and drops the now-redundant proofErr markers.
# ---------------------------------------------------------------------------
$c.Find.Execute("from the Resources folder on the course’s Moodle page", $false, $false, $false, $false, $false, $true, 1, $false, "from the Resources folder on the course’s Moodle page", 2)
$c.Find.Execute("Run the executable batch file runServer as shown below", $false, $false, $false, $false, $false, $true, 1, $false, "Run the executable batch file runServer as shown below", 2)
$c.Find.Execute("Open the HSQLDB Client Interface by running the batch file: runManagerSwing this will open the following window", $false, $false, $false, $false, $false, $true, 1, $false, "Open the HSQLDB Client Interface by running the batch file: runManagerSwing this will open the following window", 2)
$c.Find.Execute("Alter the output so that entire table contents is output properly using HTML table tags", $false, $false, $false, $false, $false, $true, 1, $false, "Alter the output so that entire table contents is output properly using HTML table tags", 2)
$c.Find.Execute(" the ResultSet and the ResultSetMetaData objects.", $false, $false, $false, $false, $false, $true, 1, $false, " the ResultSet and the ResultSetMetaData objects.", 2)

# ---------------------------------------------------------------------------
# Update the "CREATE TABLE BOOK" SQL script shown in the Week 1 instructions
# so it uses lower-case SQL keywords/identifiers and includes an explicit
# primary-key constraint on book_id (per commit message: "Update
# instructions to include syntax for PK").
# ---------------------------------------------------------------------------

# Paragraph 23: "CREATE TABLE BOOK " -> "create table Book("
$d.Paragraphs.Item(23).Range.Text = "create table Book("

# Paragraph 24 previously held just "(" ; it now becomes the book_id/PK line.
$d.Paragraphs.Item(24).Range.Text = "`tbook_id int primary key,"

# Paragraphs 25-28 previously held the BOOK_ID / TITLE / AUTHOR / GENRE
# column definitions (each already starts with a tab character) -- shift
# every column name down one row and lower-case it.
$p25 = $d.Paragraphs.Item(25)
$d.Range($p25.Range.Start + 1, $p25.Range.End).Text = "title varchar(50),"

$p26 = $d.Paragraphs.Item(26)
$d.Range($p26.Range.Start + 1, $p26.Range.End).Text = "author varchar(50),"

$p27 = $d.Paragraphs.Item(27)
$d.Range($p27.Range.Start + 1, $p27.Range.End).Text = "genre varchar(50),"

$p28 = $d.Paragraphs.Item(28)
$d.Range($p28.Range.Start + 1, $p28.Range.End).Text = "isbn varchar(50),"

# Paragraph 29 previously held "ISBN varchar(50)," ; it now becomes the
# blurb line, and the old "BLURB varchar(250)" paragraph (30) is removed
# entirely by deleting its range (merging it away).
$p29 = $d.Paragraphs.Item(29)
$d.Range($p29.Range.Start + 1, $p29.Range.End).Text = "blurb varchar(250)"

$d.Paragraphs.Item(30).Range.Delete()
